# Activities Classic - 26 June
#
# The "ActivityList" sheet's pre-set template list included an entry for
# "My Last 100 days Activities" that is being retired. Remove that whole
# row (shifting the rows below it up), tidy up the formatting that had been
# left on the "Activities Filtered by Andrew B" row, and leave the
# selection where the user left off editing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the explicit (and visually redundant) formatting that had been
# applied to the "Activities Filtered by Andrew B" cell so it goes back to
# the sheet's default/normal style.
$ws.Range("A2").Style = "Normal"

# Remove the "My Last 100 days Activities" row entirely - it's row 6 in the
# list (Pre Set Templates / Activities Filtered by Andrew B / My Activities
# Today / My Activities Tomorrow / My Activities Yesterday / My Last 100
# days Activities / ...). Deleting the whole row shifts everything below it
# up by one.
$ws.Rows(6).Delete()

# Leave the selection on B18, matching where editing left off.
[void]$ws.Range("B18").Select()
